$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.257.61'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '1.981.17'
$ws.Range('E3').Value = '  +5.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9984'
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.54'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9973'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5121'
$ws.Range('E7').Value = '  +1.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4094'
$ws.Range('E8').Value = '  +3.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08423'
$ws.Range('E9').Value = '  +2.70%  '
$ws.Range('E10').Value = '  +3.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.53'
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.27'
$ws.Range('E12').Value = '  +3.63%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.957.24'
$ws.Range('E13').Value = '  +4.24%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.461'
$ws.Range('E14').Value = '  +2.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.392'
$ws.Range('E15').Value = '  +2.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.75'
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001106'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06528'
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.73'
$ws.Range('E20').Value = '  +3.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9977'
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.068'
$ws.Range('E22').Value = '  +3.70%  '
$ws.Range('D23').Value = '30.312.43'
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('E24').Value = '  +2.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.195'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.183.16'
$ws.Range('E26').Value = '  +4.38%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.55'
$ws.Range('E27').Value = '  +6.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.62'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.370'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '130.26'
$ws.Range('E30').Value = '  +2.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.140'
$ws.Range('E31').Value = '  +6.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1051'
$ws.Range('E32').Value = '  +1.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.034'
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.792'
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.305'
$ws.Range('E35').Value = '  +11.22%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02474'
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.367'
$ws.Range('E37').Value = '  +2.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06493'
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.897'
$ws.Range('E40').Value = '  +4.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6570'
$ws.Range('E41').Value = '  +4.19%  '
$ws.Range('E42').Value = '  +4.04%  '
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.57'
$ws.Range('E44').Value = '  +4.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6112'
$ws.Range('E45').Value = '  +3.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.188'
$ws.Range('E46').Value = '  +4.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.637'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.222'
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.48'
$ws.Range('E49').Value = '  +0.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '79.45'
$ws.Range('E50').Value = '  +2.54%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06900'
$ws.Range('E51').Value = '  +2.07%  '
